$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 6 with trade data (values first)
$ws.Cells.Item(6, 1).Value = 42647.681828703702
$ws.Cells.Item(6, 2).Value = $false
$ws.Cells.Item(6, 3).Value = 9971.89
$ws.Cells.Item(6, 4).Value = 10013.450000000001
$ws.Cells.Item(6, 5).Value = 18.12
$ws.Cells.Item(6, 6).Value = 17.97
$ws.Cells.Item(6, 7).Value = $false
$ws.Cells.Item(6, 8).Value = -0.83
$ws.Cells.Item(6, 9).Value = $false

# Copy formatting (style s="1", date format) from row 5's matching cells
$ws.Cells.Item(5, 1).Copy()
$ws.Cells.Item(6, 1).PasteSpecial(-4122)

$ws.Cells.Item(5, 7).Copy()
$ws.Cells.Item(6, 7).PasteSpecial(-4122)
